$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 17 and row 18
$ws.Range("B17").Value = 6815306
$ws.Range("B18").Value = 6815303
$ws.Range("F17").Value = 'OFK Petrovac'
$ws.Range("F18").Value = 'FK Mornar Bar'
$ws.Range("G17").Value = 'FK Rudar Pljevlja'
$ws.Range("G18").Value = 'FK Arsenal'
$ws.Range("H17").Value = 2
$ws.Range("H18").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("I18").Value = 1
$ws.Range("J17").Value = 'H'
$ws.Range("J18").Value = 'A'
$ws.Range("K17").Value = 2.1
$ws.Range("K18").Value = 2.4
$ws.Range("L17").Value = 3.1
$ws.Range("L18").Value = 3
$ws.Range("M17").Value = 3.2
$ws.Range("M18").Value = 2.75
$ws.Range("N17").Value = 1.615
$ws.Range("N18").Value = 2.4
$ws.Range("O17").Value = 3.5
$ws.Range("O18").Value = 3
$ws.Range("P17").Value = 5
$ws.Range("P18").Value = 2.75
$ws.Range("Q17").Value = -0.75
$ws.Range("Q18").Value = 0
$ws.Range("R17").Value = 1.85
$ws.Range("R18").Value = 1.775
$ws.Range("S17").Value = 1.95
$ws.Range("S18").Value = 2.025
$ws.Range("T17").Value = 2.25
$ws.Range("T18").Value = 2
$ws.Range("W17").Value = 0.615
$ws.Range("W18").Value = -1
$ws.Range("Y17").Value = -1
$ws.Range("Y18").Value = 1.75
$ws.Range("Z17").Value = 0.8500000000000001
$ws.Range("Z18").Value = -1
$ws.Range("AA17").Value = -1
$ws.Range("AA18").Value = 1.025
$ws.Range("AB17").Value = -0.5
$ws.Range("AB18").Value = -1
$ws.Range("AC17").Value = 0.475
$ws.Range("AC18").Value = 0.95

# Swap row 20 and row 21
$ws.Range("B20").Value = 6815305
$ws.Range("B21").Value = 6815304
$ws.Range("F20").Value = 'Buducnost Podgorica'
$ws.Range("F21").Value = 'FK Jedinstvo Bijelo Polje'
$ws.Range("G20").Value = 'FK Jezero'
$ws.Range("G21").Value = 'Sutjeska Niksic'
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 0
$ws.Range("I20").Value = 1
$ws.Range("I21").Value = 0
$ws.Range("K20").Value = 1.3
$ws.Range("K21").Value = 5.5
$ws.Range("L20").Value = 5
$ws.Range("L21").Value = 3.75
$ws.Range("M20").Value = 7
$ws.Range("M21").Value = 1.5
$ws.Range("N20").Value = 1.571
$ws.Range("N21").Value = 3.6
$ws.Range("O20").Value = 4
$ws.Range("O21").Value = 3.2
$ws.Range("P20").Value = 4.2
$ws.Range("P21").Value = 1.909
$ws.Range("Q20").Value = -0.75
$ws.Range("Q21").Value = 0.5
$ws.Range("R20").Value = 1.75
$ws.Range("R21").Value = 1.825
$ws.Range("S20").Value = 1.95
$ws.Range("S21").Value = 1.975
$ws.Range("T20").Value = 2.5
$ws.Range("T21").Value = 2.25
$ws.Range("U20").Value = 1.95
$ws.Range("U21").Value = 1.875
$ws.Range("V20").Value = 1.85
$ws.Range("V21").Value = 1.925
$ws.Range("X20").Value = 3
$ws.Range("X21").Value = 2.2
$ws.Range("Z20").Value = -1
$ws.Range("Z21").Value = 0.825
$ws.Range("AA20").Value = 0.95
$ws.Range("AA21").Value = -1
$ws.Range("AC20").Value = 0.8500000000000001
$ws.Range("AC21").Value = 0.925

# Swap row 59 and row 60
$ws.Range("B59").Value = 6815338
$ws.Range("B60").Value = 6815427
$ws.Range("F59").Value = 'OFK Petrovac'
$ws.Range("F60").Value = 'FK Mornar Bar'
$ws.Range("G59").Value = 'FK Decic Tuzi'
$ws.Range("G60").Value = 'OFK Mladost DG'
$ws.Range("I59").Value = 3
$ws.Range("I60").Value = 1
$ws.Range("J59").Value = 'A'
$ws.Range("J60").Value = 'H'
$ws.Range("K59").Value = 2.625
$ws.Range("K60").Value = 1.833
$ws.Range("L59").Value = 2.875
$ws.Range("L60").Value = 3.1
$ws.Range("M59").Value = 2.6
$ws.Range("M60").Value = 4
$ws.Range("N59").Value = 3.1
$ws.Range("N60").Value = 1.833
$ws.Range("O59").Value = 2.9
$ws.Range("O60").Value = 3.1
$ws.Range("P59").Value = 2.25
$ws.Range("P60").Value = 4
$ws.Range("Q59").Value = 0.25
$ws.Range("Q60").Value = -0.5
$ws.Range("R59").Value = 1.8
$ws.Range("R60").Value = 1.875
$ws.Range("S59").Value = 2
$ws.Range("S60").Value = 1.925
$ws.Range("T59").Value = 2.25
$ws.Range("T60").Value = 2
$ws.Range("U59").Value = 1.975
$ws.Range("U60").Value = 1.775
$ws.Range("V59").Value = 1.725
$ws.Range("V60").Value = 2.025
$ws.Range("W59").Value = -1
$ws.Range("W60").Value = 0.833
$ws.Range("Y59").Value = 1.25
$ws.Range("Y60").Value = -1
$ws.Range("Z59").Value = -1
$ws.Range("Z60").Value = 0.875
$ws.Range("AA59").Value = 1
$ws.Range("AA60").Value = -1
$ws.Range("AB59").Value = 0.9750000000000001
$ws.Range("AB60").Value = 0.7749999999999999

# Swap row 62 and row 63
$ws.Range("B62").Value = 6815343
$ws.Range("B63").Value = 7366683
$ws.Range("F62").Value = 'Sutjeska Niksic'
$ws.Range("F63").Value = 'FK Arsenal'
$ws.Range("G62").Value = 'FK Jedinstvo Bijelo Polje'
$ws.Range("G63").Value = 'FK Mornar Bar'
$ws.Range("I62").Value = 0
$ws.Range("I63").Value = 2
$ws.Range("J62").Value = 'H'
$ws.Range("J63").Value = 'D'
$ws.Range("K62").Value = 1.333
$ws.Range("K63").Value = 2.375
$ws.Range("L62").Value = 4.2
$ws.Range("L63").Value = 2.8
$ws.Range("M62").Value = 8
$ws.Range("M63").Value = 3
$ws.Range("N62").Value = 1.333
$ws.Range("N63").Value = 2.3
$ws.Range("O62").Value = 4.2
$ws.Range("O63").Value = 2.7
$ws.Range("P62").Value = 8
$ws.Range("P63").Value = 3.3
$ws.Range("Q62").Value = -1.5
$ws.Range("Q63").Value = -0.25
$ws.Range("R62").Value = 1.975
$ws.Range("R63").Value = 2
$ws.Range("S62").Value = 1.825
$ws.Range("S63").Value = 1.8
$ws.Range("T62").Value = 2.75
$ws.Range("T63").Value = 1.75
$ws.Range("U62").Value = 1.9
$ws.Range("U63").Value = 1.875
$ws.Range("V62").Value = 1.9
$ws.Range("V63").Value = 1.925
$ws.Range("W62").Value = 0.333
$ws.Range("W63").Value = -1
$ws.Range("X62").Value = -1
$ws.Range("X63").Value = 1.7
$ws.Range("Z62").Value = 0.9750000000000001
$ws.Range("Z63").Value = -0.5
$ws.Range("AA62").Value = -1
$ws.Range("AA63").Value = 0.4
$ws.Range("AB62").Value = -1
$ws.Range("AB63").Value = 0.875
$ws.Range("AC62").Value = 0.8999999999999999
$ws.Range("AC63").Value = -1

# Swap row 81 and row 82
$ws.Range("B81").Value = 6815362
$ws.Range("B82").Value = 6815430
$ws.Range("F81").Value = 'Sutjeska Niksic'
$ws.Range("F82").Value = 'Buducnost Podgorica'
$ws.Range("G81").Value = 'FK Decic Tuzi'
$ws.Range("G82").Value = 'FK Mornar Bar'
$ws.Range("H81").Value = 1
$ws.Range("H82").Value = 4
$ws.Range("I81").Value = 1
$ws.Range("I82").Value = 3
$ws.Range("J81").Value = 'D'
$ws.Range("J82").Value = 'H'
$ws.Range("K81").Value = 2.2
$ws.Range("K82").Value = 1.444
$ws.Range("L81").Value = 3
$ws.Range("L82").Value = 3.75
$ws.Range("M81").Value = 3.1
$ws.Range("M82").Value = 6.5
$ws.Range("N81").Value = 2.375
$ws.Range("N82").Value = 1.4
$ws.Range("O81").Value = 2.875
$ws.Range("O82").Value = 4
$ws.Range("P81").Value = 3
$ws.Range("P82").Value = 7
$ws.Range("Q81").Value = -0.25
$ws.Range("Q82").Value = -1.25
$ws.Range("R81").Value = 2.05
$ws.Range("R82").Value = 1.875
$ws.Range("S81").Value = 1.75
$ws.Range("S82").Value = 1.925
$ws.Range("T81").Value = 2
$ws.Range("T82").Value = 2.5
$ws.Range("U81").Value = 1.8
$ws.Range("U82").Value = 1.775
$ws.Range("V81").Value = 2
$ws.Range("V82").Value = 1.925
$ws.Range("W81").Value = -1
$ws.Range("W82").Value = 0.3999999999999999
$ws.Range("X81").Value = 1.875
$ws.Range("X82").Value = -1
$ws.Range("AA81").Value = 0.375
$ws.Range("AA82").Value = 0.4625
$ws.Range("AB81").Value = 0
$ws.Range("AB82").Value = 0.7749999999999999
$ws.Range("AC81").Value = -0
$ws.Range("AC82").Value = -1

# Swap row 107 and row 108
$ws.Range("B107").Value = 7890506
$ws.Range("B108").Value = 7890508
$ws.Range("F107").Value = 'FK Mornar Bar'
$ws.Range("F108").Value = 'OFK Petrovac'
$ws.Range("G107").Value = 'FK Arsenal'
$ws.Range("G108").Value = 'FK Rudar Pljevlja'
$ws.Range("H107").Value = 0
$ws.Range("H108").Value = 1
$ws.Range("I107").Value = 0
$ws.Range("I108").Value = 1
$ws.Range("K107").Value = 1.85
$ws.Range("K108").Value = 1.75
$ws.Range("M107").Value = 3.9
$ws.Range("M108").Value = 4.5
$ws.Range("N107").Value = 1.85
$ws.Range("N108").Value = 1.8
$ws.Range("O107").Value = 3.3
$ws.Range("O108").Value = 3.2
$ws.Range("P107").Value = 3.5
$ws.Range("P108").Value = 4
$ws.Range("R107").Value = 1.925
$ws.Range("R108").Value = 1.875
$ws.Range("S107").Value = 1.875
$ws.Range("S108").Value = 1.925
$ws.Range("T107").Value = 2
$ws.Range("T108").Value = 2.25
$ws.Range("X107").Value = 2.3
$ws.Range("X108").Value = 2.2
$ws.Range("AA107").Value = 0.875
$ws.Range("AA108").Value = 0.925
$ws.Range("AB107").Value = -1
$ws.Range("AB108").Value = -0.5
$ws.Range("AC107").Value = 0.8500000000000001
$ws.Range("AC108").Value = 0.425

# Swap row 121 and row 122
$ws.Range("B121").Value = 6815397
$ws.Range("B122").Value = 6815398
$ws.Range("F121").Value = 'FK Arsenal'
$ws.Range("F122").Value = 'FK Jedinstvo Bijelo Polje'
$ws.Range("G121").Value = 'OFK Petrovac'
$ws.Range("G122").Value = 'Buducnost Podgorica'
$ws.Range("I121").Value = 1
$ws.Range("I122").Value = 2
$ws.Range("J121").Value = 'D'
$ws.Range("J122").Value = 'A'
$ws.Range("K121").Value = 2.6
$ws.Range("K122").Value = 5.75
$ws.Range("L121").Value = 2.7
$ws.Range("L122").Value = 4
$ws.Range("M121").Value = 2.8
$ws.Range("M122").Value = 1.444
$ws.Range("N121").Value = 2.75
$ws.Range("N122").Value = 6.5
$ws.Range("O121").Value = 2.5
$ws.Range("O122").Value = 4.2
$ws.Range("P121").Value = 2.875
$ws.Range("P122").Value = 1.4
$ws.Range("Q121").Value = 0
$ws.Range("Q122").Value = 1.25
$ws.Range("R121").Value = 1.85
$ws.Range("R122").Value = 1.9
$ws.Range("S121").Value = 1.95
$ws.Range("S122").Value = 1.9
$ws.Range("T121").Value = 2
$ws.Range("T122").Value = 2.75
$ws.Range("U121").Value = 2.025
$ws.Range("U122").Value = 1.875
$ws.Range("V121").Value = 1.775
$ws.Range("V122").Value = 1.925
$ws.Range("X121").Value = 1.5
$ws.Range("X122").Value = -1
$ws.Range("Y121").Value = -1
$ws.Range("Y122").Value = 0.3999999999999999
$ws.Range("Z121").Value = 0
$ws.Range("Z122").Value = 0.45
$ws.Range("AA121").Value = -0
$ws.Range("AA122").Value = -0.5
$ws.Range("AB121").Value = 0
$ws.Range("AB122").Value = 0.4375
$ws.Range("AC121").Value = -0
$ws.Range("AC122").Value = -0.5

# Swap row 126 and row 127
$ws.Range("B126").Value = 6815402
$ws.Range("B127").Value = 6815401
$ws.Range("F126").Value = 'FK Rudar Pljevlja'
$ws.Range("F127").Value = 'FK Decic Tuzi'
$ws.Range("G126").Value = 'FK Jezero'
$ws.Range("G127").Value = 'Sutjeska Niksic'
$ws.Range("I126").Value = 1
$ws.Range("I127").Value = 0
$ws.Range("J126").Value = 'A'
$ws.Range("J127").Value = 'D'
$ws.Range("K126").Value = 2.8
$ws.Range("K127").Value = 2.55
$ws.Range("M126").Value = 2.375
$ws.Range("M127").Value = 2.6
$ws.Range("N126").Value = 2.45
$ws.Range("N127").Value = 2.1
$ws.Range("O126").Value = 2.9
$ws.Range("O127").Value = 3.1
$ws.Range("P126").Value = 2.75
$ws.Range("P127").Value = 3.3
$ws.Range("Q126").Value = 0
$ws.Range("Q127").Value = -0.25
$ws.Range("R126").Value = 1.775
$ws.Range("R127").Value = 1.825
$ws.Range("S126").Value = 2.025
$ws.Range("S127").Value = 1.975
$ws.Range("T126").Value = 1.75
$ws.Range("T127").Value = 2
$ws.Range("U126").Value = 1.825
$ws.Range("U127").Value = 1.925
$ws.Range("V126").Value = 1.975
$ws.Range("V127").Value = 1.875
$ws.Range("X126").Value = -1
$ws.Range("X127").Value = 2.1
$ws.Range("Y126").Value = 1.75
$ws.Range("Y127").Value = -1
$ws.Range("Z126").Value = -1
$ws.Range("Z127").Value = -0.5
$ws.Range("AA126").Value = 1.025
$ws.Range("AA127").Value = 0.4875
$ws.Range("AC126").Value = 0.9750000000000001
$ws.Range("AC127").Value = 0.875

